# Update the StructureDefinition workbook to point at the new
# LinuxForHealth home (was: IBM/Alvearie), and bump the version/date
# metadata to match the new published IG build.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet -------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/default-value"
# Version
$wsMeta.Range("B3").Value = "8.0.0"
# Date
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet ---------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# The Extension.url row's "Fixed Value" column mirrors the URL above.
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/default-value"

# The root "Extension" row no longer carries the ele-1/ext-1 constraint
# text in its "Constraint(s)" column (it now only shows up on the child
# rows, e.g. Extension.extension / Extension.value[x]).
$wsElem.Range("AI2").Value = ""
